{"js": "// Helper: replace the first search hit's text (optionally set bold) leaving\n// other run formatting untouched. Operates on the matched Range only, so it\n// never touches the paragraph mark's rPr.\nasync function replaceFirst(body, searchText, newText) {\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nasync function deleteParagraphContaining(body, searchText) {\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  const paras = results.items[0].paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n  paras.items[0].delete();\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1) Title font size 26 -> 24 half-points (13pt -> 12pt), run-scoped so the\n//    paragraph mark is left alone.\n{\n  const results = body.search(\"INFORME FINAL DE AUDITOR\u00cdA INFORM\u00c1TICA\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].font.size = 12;\n  await context.sync();\n}\n\n// 2) Subtitle line: shorten department name + append the audit date.\nawait replaceFirst(\n  body,\n  \"Departamento de Infraestructura de Red y Comunicaciones - TeleNET Solutions C.A.\",\n  \"Infraestructura de Red - TeleNET Solutions C.A. | 03-07 febrero 2026\"\n);\n\n// 3) Drop the \"Fecha de auditor\u00eda: 03 al 07 de febrero de 2026 | \" lead-in,\n//    keeping \"Equipo Auditor: \" (bold) + the auditor list that follows.\nawait replaceFirst(body, \"Fecha de auditor\u00eda: 03 al 07 de febrero de 2026 | \", \"\");\n\n// Abbreviate the auditor names.\nawait replaceFirst(\n  body,\n  \"V\u00edctor Ysea (L\u00edder), Mar\u00eda Yoris, Fiorella David, Juan Marcano\",\n  \"V. Ysea (L\u00edder), M. Yoris, F. David, J. Marcano\"\n);\n\n// 4) Remove the whole \"Personal Entrevistado: ...\" paragraph.\nawait deleteParagraphContaining(body, \"Personal Entrevistado: \");\n\n// 5) Collapse \"ALCANCE Y OBJETIVOS\" heading + body paragraph into a single\n//    paragraph that starts with a bold \"ALCANCE: \" label.\nawait deleteParagraphContaining(body, \"ALCANCE Y OBJETIVOS\");\n{\n  const results = body.search(\"Se evaluaron los controles de seguridad\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  const paras = results.items[0].paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n  const para = paras.items[0];\n\n  para.getRange().insertText(\n    \"Evaluaci\u00f3n de seguridad, infraestructura, controles internos, cumplimiento y eficiencia operativa. Muestra: semana 06-10 enero 2026.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n\n  const label = para.insertText(\"ALCANCE: \", Word.InsertLocation.start);\n  label.font.bold = true;\n  await context.sync();\n}\n\n// 6) Collapse \"HALLAZGOS PRINCIPALES\" heading + body paragraph.\nawait deleteParagraphContaining(body, \"HALLAZGOS PRINCIPALES\");\n{\n  const results = body.search(\"Se identificaron seis hallazgos significativos\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  const paras = results.items[0].paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n  const para = paras.items[0];\n\n  para.getRange().insertText(\n    \"(1) Pruebas de respaldo sin ejecutar en 8 meses; (2) Licencia Nessus vencida 3 meses; (3) Firmware desactualizado 2+ a\u00f1os; (4) UPS sin mantenimiento 18 meses; (5) Cambios sin aprobaci\u00f3n formal; (6) Servidores HP sin soporte.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n\n  const label = para.insertText(\"HALLAZGOS: \", Word.InsertLocation.start);\n  label.font.bold = true;\n  await context.sync();\n}\n\n// 7) Collapse \"CONCLUSIONES\" heading + body paragraph.\nawait deleteParagraphContaining(body, \"CONCLUSIONES\");\n{\n  const results = body.search(\"El \u00e1rea presenta deficiencias significativas\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  const paras = results.items[0].paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n  const para = paras.items[0];\n\n  para.getRange().insertText(\n    \"\u00c1rea con deficiencias significativas. 3/5 \u00e1reas DEFICIENTES, 2/5 PARCIALES. Dos recomendaciones de auditor\u00eda 2024 sin implementar.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n\n  const label = para.insertText(\"CONCLUSI\u00d3N: \", Word.InsertLocation.start);\n  label.font.bold = true;\n  await context.sync();\n}\n\n// 8) Collapse \"RECOMENDACIONES\" heading + body paragraph.\nawait deleteParagraphContaining(body, \"RECOMENDACIONES\");\n{\n  const results = body.search(\"Se recomienda: (1) implementar pruebas mensuales\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  const paras = results.items[0].paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n  const para = paras.items[0];\n\n  para.getRange().insertText(\n    \"(1) Pruebas de respaldo mensuales - inmediato; (2) Renovar Nessus - 15 d\u00edas; (3) Plan actualizaci\u00f3n firmware - 30 d\u00edas; (4) Mantenimiento UPS - inmediato; (5) Cumplir gesti\u00f3n de cambios - inmediato; (6) Plan migraci\u00f3n servidores - 60 d\u00edas. Seguimiento: 90 d\u00edas.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n\n  const label = para.insertText(\"RECOMENDACIONES: \", Word.InsertLocation.start);\n  label.font.bold = true;\n  await context.sync();\n}\n\n// 9) Shorten the signature underscores.\nawait replaceFirst(\n  body,\n  \"____________________________                    ____________________________\",\n  \"_______________________                    _______________________\"\n);\n\n// 10) Tighten the spacing between the two signer names.\nawait replaceFirst(\n  body,\n  \"V\u00edctor Ysea                                                    Carlos Mendoza\",\n  \"V\u00edctor Ysea                                        Carlos Mendoza\"\n);\n\n// 11) Shrink the signature caption font (20 -> 18 half-points) and shorten\n//     the caption text, run-scoped.\n{\n  const results = body.search(\"L\u00edder del Equipo Auditor                                Gerente de Infraestructura\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  const r = results.items[0];\n  r.font.size = 9;\n  r.insertText(\"L\u00edder Equipo Auditor                          Gerente Infraestructura\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Find-Replace($doc, $findText, $replaceText) {\n    $range = $doc.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\nfunction Get-FirstMatchRange($doc, $findText) {\n    $range = $doc.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n    return $range\n}\n\nfunction Delete-ParagraphContaining($doc, $findText) {\n    $range = Get-FirstMatchRange $doc $findText\n    # `Find` only returns the matched span; Expand to the enclosing paragraph\n    # (wdParagraph = 4) before deleting so the whole line (and its mark) goes.\n    $range.Expand(4)\n    $range.Delete()\n}\n\nfunction Set-LeadingBoldLabel($doc, $bodyText, $label) {\n    $range = Get-FirstMatchRange $doc $bodyText\n    $range.Expand(4)\n    $insertRange = $range.Duplicate()\n    $insertRange.Collapse(1)\n    $insertRange.InsertBefore($label)\n    $insertRange.Font.Bold = 1\n}\n\n# 1) Title font size 26 -> 24 half-points (13pt -> 12pt); Find-based range so\n#    only the run is touched, not the paragraph mark.\n$titleRange = Get-FirstMatchRange $d \"INFORME FINAL DE AUDITOR\u00cdA INFORM\u00c1TICA\"\n$titleRange.Font.Size = 12\n\n# 2) Subtitle line: shorten department name + append the audit date.\nFind-Replace $d \"Departamento de Infraestructura de Red y Comunicaciones - TeleNET Solutions C.A.\" \"Infraestructura de Red - TeleNET Solutions C.A. | 03-07 febrero 2026\"\n\n# 3) Drop the \"Fecha de auditor\u00eda: 03 al 07 de febrero de 2026 | \" lead-in,\n#    keeping \"Equipo Auditor: \" (bold) + the auditor list that follows.\nFind-Replace $d \"Fecha de auditor\u00eda: 03 al 07 de febrero de 2026 | \" \"\"\n\n# Abbreviate the auditor names.\nFind-Replace $d \"V\u00edctor Ysea (L\u00edder), Mar\u00eda Yoris, Fiorella David, Juan Marcano\" \"V. Ysea (L\u00edder), M. Yoris, F. David, J. Marcano\"\n\n# 4) Remove the whole \"Personal Entrevistado: ...\" paragraph.\nDelete-ParagraphContaining $d \"Personal Entrevistado: \"\n\n# 5) Collapse \"ALCANCE Y OBJETIVOS\" heading + body paragraph into a single\n#    paragraph that starts with a bold \"ALCANCE: \" label.\nDelete-ParagraphContaining $d \"ALCANCE Y OBJETIVOS\"\nFind-Replace $d \"Se evaluaron los controles de seguridad, infraestructura tecnol\u00f3gica, controles internos, cumplimiento normativo y eficiencia operativa del \u00e1rea, utilizando como muestra la semana del 06 al 10 de enero de 2026. El objetivo fue verificar que los recursos inform\u00e1ticos son adecuadamente utilizados y vigilados, identificando vulnerabilidades y evaluando la capacidad de continuidad operativa.\" \"Evaluaci\u00f3n de seguridad, infraestructura, controles internos, cumplimiento y eficiencia operativa. Muestra: semana 06-10 enero 2026.\"\nSet-LeadingBoldLabel $d \"Evaluaci\u00f3n de seguridad, infraestructura, controles internos, cumplimiento y eficiencia operativa. Muestra: semana 06-10 enero 2026.\" \"ALCANCE: \"\n\n# 6) Collapse \"HALLAZGOS PRINCIPALES\" heading + body paragraph.\nDelete-ParagraphContaining $d \"HALLAZGOS PRINCIPALES\"\nFind-Replace $d \"Se identificaron seis hallazgos significativos: (1) ausencia de pruebas de restauraci\u00f3n de respaldos desde hace 8 meses, incumpliendo la pol\u00edtica de pruebas mensuales; (2) licencia de Nessus vencida hace 3 meses, inhabilitando el escaneo de vulnerabilidades; (3) firmware de equipos Cisco desactualizado por m\u00e1s de 2 a\u00f1os con vulnerabilidades conocidas; (4) una unidad UPS sin mantenimiento preventivo en 18 meses; (5) tres cambios ejecutados sin aprobaci\u00f3n formal durante la semana de muestreo; y (6) dos servidores HP operando sin soporte del fabricante.\" \"(1) Pruebas de respaldo sin ejecutar en 8 meses; (2) Licencia Nessus vencida 3 meses; (3) Firmware desactualizado 2+ a\u00f1os; (4) UPS sin mantenimiento 18 meses; (5) Cambios sin aprobaci\u00f3n formal; (6) Servidores HP sin soporte.\"\nSet-LeadingBoldLabel $d \"(1) Pruebas de respaldo sin ejecutar en 8 meses; (2) Licencia Nessus vencida 3 meses; (3) Firmware desactualizado 2+ a\u00f1os; (4) UPS sin mantenimiento 18 meses; (5) Cambios sin aprobaci\u00f3n formal; (6) Servidores HP sin soporte.\" \"HALLAZGOS: \"\n\n# 7) Collapse \"CONCLUSIONES\" heading + body paragraph.\nDelete-ParagraphContaining $d \"CONCLUSIONES\"\nFind-Replace $d \"El \u00e1rea presenta deficiencias significativas: tres de cinco \u00e1reas evaluadas resultaron DEFICIENTES (Seguridad, Infraestructura y Controles Internos) y dos PARCIALES (Cumplimiento y Eficiencia). Se evidenci\u00f3 adem\u00e1s que dos recomendaciones de la auditor\u00eda de 2024 no fueron implementadas.\" \"\u00c1rea con deficiencias significativas. 3/5 \u00e1reas DEFICIENTES, 2/5 PARCIALES. Dos recomendaciones de auditor\u00eda 2024 sin implementar.\"\nSet-LeadingBoldLabel $d \"\u00c1rea con deficiencias significativas. 3/5 \u00e1reas DEFICIENTES, 2/5 PARCIALES. Dos recomendaciones de auditor\u00eda 2024 sin implementar.\" \"CONCLUSI\u00d3N: \"\n\n# 8) Collapse \"RECOMENDACIONES\" heading + body paragraph.\nDelete-ParagraphContaining $d \"RECOMENDACIONES\"\nFind-Replace $d \"Se recomienda: (1) implementar pruebas mensuales de restauraci\u00f3n de respaldos de manera inmediata; (2) renovar la licencia de Nessus o implementar alternativa en 15 d\u00edas; (3) elaborar plan de actualizaci\u00f3n de firmware en 30 d\u00edas; (4) ejecutar mantenimiento de UPS de forma inmediata; (5) reforzar el cumplimiento del procedimiento de gesti\u00f3n de cambios; y (6) presentar plan de migraci\u00f3n de servidores obsoletos en 60 d\u00edas. El seguimiento se realizar\u00e1 en 90 d\u00edas.\" \"(1) Pruebas de respaldo mensuales - inmediato; (2) Renovar Nessus - 15 d\u00edas; (3) Plan actualizaci\u00f3n firmware - 30 d\u00edas; (4) Mantenimiento UPS - inmediato; (5) Cumplir gesti\u00f3n de cambios - inmediato; (6) Plan migraci\u00f3n servidores - 60 d\u00edas. Seguimiento: 90 d\u00edas.\"\nSet-LeadingBoldLabel $d \"(1) Pruebas de respaldo mensuales - inmediato; (2) Renovar Nessus - 15 d\u00edas; (3) Plan actualizaci\u00f3n firmware - 30 d\u00edas; (4) Mantenimiento UPS - inmediato; (5) Cumplir gesti\u00f3n de cambios - inmediato; (6) Plan migraci\u00f3n servidores - 60 d\u00edas. Seguimiento: 90 d\u00edas.\" \"RECOMENDACIONES: \"\n\n# 9) Shorten the signature underscores.\nFind-Replace $d \"____________________________                    ____________________________\" \"_______________________                    _______________________\"\n\n# 10) Tighten the spacing between the two signer names.\nFind-Replace $d \"V\u00edctor Ysea                                                    Carlos Mendoza\" \"V\u00edctor Ysea                                        Carlos Mendoza\"\n\n# 11) Shrink the signature caption font (20 -> 18 half-points) and shorten\n#     the caption text, run-scoped via Find.\n$captionRange = Get-FirstMatchRange $d \"L\u00edder del Equipo Auditor                                Gerente de Infraestructura\"\n$captionRange.Font.Size = 9\nFind-Replace $d \"L\u00edder del Equipo Auditor                                Gerente de Infraestructura\" \"L\u00edder Equipo Auditor                          Gerente Infraestructura\"\n"}
